# Fixing memory leaks in string functions
#
# - Remove the now-obsolete "Bank Segment" / "Segments Per Bank" / "Segment Calc"
#   rows from the Dynamic sheet.
# - Replace the broken (#REF!) allocation rows 3-5 on Bank3 with a simple "TBD"
#   placeholder in column B, clearing the now-unused C:E columns.
# - Bank3 becomes the active/selected sheet & cell selections are updated.

$wb = $excel.ActiveWorkbook

# --- Dynamic sheet: drop rows 14-16 (Bank Segment / Segments Per Bank / Segment Calc) ---
$wsDynamic = $wb.Worksheets.Item("Dynamic")
$wsDynamic.Rows("14:16").Delete()
$wsDynamic.Range("B13").Select() | Out-Null

# --- Bank3 sheet: rows 3-5 collapse to a "TBD" placeholder ---
$wsBank3 = $wb.Worksheets.Item("Bank3")
$wsBank3.Range("B3").Value = "TBD"
$wsBank3.Range("C3:E3").ClearContents()
$wsBank3.Range("B4").Value = "TBD"
$wsBank3.Range("C4:E4").ClearContents()
$wsBank3.Range("B5").Value = "TBD"
$wsBank3.Range("C5:E5").ClearContents()

# Bank3 is the newly-active tab, with D5 selected
$wsBank3.Activate() | Out-Null
$wsBank3.Range("D5").Select() | Out-Null
